# Update the "Latest Handoff Datetime" column for each localized language
# sheet to reflect the new handoff report generated on 2016-02-23.
#
# Each language worksheet has rows 2-59 in column D all sharing the same
# "Latest Handoff Datetime" text value. We replace the old datetime string
# with the new one for every sheet below.

$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Sheet = "ja-jp"; NewValue = "2016-02-23 03:01:00" },
    @{ Sheet = "de-de"; NewValue = "2016-02-23 03:01:20" },
    @{ Sheet = "fr-fr"; NewValue = "2016-02-23 03:01:38" },
    @{ Sheet = "zh-cn"; NewValue = "2016-02-23 03:01:57" },
    @{ Sheet = "zh-tw"; NewValue = "2016-02-23 03:02:16" },
    @{ Sheet = "ko-kr"; NewValue = "2016-02-23 03:02:35" },
    @{ Sheet = "es-es"; NewValue = "2016-02-23 03:02:55" },
    @{ Sheet = "it-it"; NewValue = "2016-02-23 03:03:15" },
    @{ Sheet = "ru-ru"; NewValue = "2016-02-23 03:03:35" },
    @{ Sheet = "pt-br"; NewValue = "2016-02-23 03:03:56" }
)

foreach ($update in $updates) {
    $ws = $wb.Worksheets.Item($update.Sheet)
    $rng = $ws.Range("D2:D59")
    $rng.Value2 = $update.NewValue
}
